$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume snapshot values.
# Columns D (Price) and E (Volume(1h)) hold numeric-looking text, so each
# target cell is switched to Text number format before assignment; this
# prevents Excel from auto-coercing strings like "244.79" or "-0.95%"
# into numeric/percentage values (and losing precision/formatting).

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '244.79'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.95%'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '27.12'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '2.73%'

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.069'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.33%'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05690'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '1.33%'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.475'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.63%'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8209'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.91%'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8397'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-0.78%'

# Row 9
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1326'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-1.10%'

# Row 10
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06904'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-1.44%'

# Row 11
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02860'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.81%'

# Row 12
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09399'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.15%'

# Row 13
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001512'

# Row 14
$ws.Range('B14').Value = 'CoinExToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.04093'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-12.33%'

# Row 15
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0005975'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-93.97%'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006094'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.96%'

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.510'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-2.18%'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.000'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.55%'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.314'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '9.24%'

# Row 20
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-1.56%'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.03189'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.92%'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.1291'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-2.24%'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.555'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-5.00%'

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.73%'

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-2.59%'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.003950'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-14.16%'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.00009795'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '2.05%'

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '-0.05%'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03706'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.78%'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.005897'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-5.33%'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1055'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.26%'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002356'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-5.73%'

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009371'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '5.24%'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005206'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-1.60%'

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.04%'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1014'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-32.32%'

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '2.35%'

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.04%'

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.04%'
